# Applies the tracked changes to "Pemrograman Web dan Mobile.docx".
$d = $word.ActiveDocument

# Searches for $searchText starting at character position $fromPos (to the end
# of the document) and returns the matched Range. wdFindContinue=1.
function Find-From($searchText, $fromPos) {
    $sub = $d.Range($fromPos, $d.Content.End)
    $ok = $sub.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Could not find text: $searchText"
    }
    return $sub
}

# --- Change 1: remove the empty paragraph right after the first table,
#     just before the "Saya Sudah Belajar dan Mengerti dan Saya Bisa" heading.
$r = Find-From "Saya Sudah Belajar dan Mengerti" 0
$emptyPara = $d.Range($r.Start - 1, $r.Start - 1)
$emptyPara.Expand(4) | Out-Null
$emptyPara.Delete()
$pos = $r.End

# --- Change 2: add a new item 7 after "6. Menginstal Prettier agar codingan terlihat rapi."
$r = Find-From "6. Menginstal Prettier agar codingan terlihat rapi." $pos
$insPt = $d.Range($r.End, $r.End)
$insPt.InsertParagraphAfter()
$newPara = $d.Range($r.End + 1, $r.End + 1)
$newPara.InsertAfter("7. Menginstal Bracket Pair agar bisa melihat pasangan dari setiap kurung yang digunakan.")
$pos = $newPara.End

# --- Change 3: fill in the first "Saya Belum Mengerti" list (items 1-3).
# Item "1." currently has a trailing "  " run followed by the _GoBack bookmark;
# replace the trailing run text and drop the bookmark from here.
$r = Find-From "Saya Belum Mengerti" $pos
$pos = $r.End

$r1 = Find-From "1." $pos
$after1 = $d.Range($r1.End, $r1.End + 2)
if ($after1.Text -eq "  ") {
    $after1.Text = " Menggunakan Boostrap."
}
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()
$pos = $after1.End

$r2 = Find-From "2." $pos
$after2 = $d.Range($r2.End, $r2.End + 2)
if ($after2.Text -eq "  ") {
    $after2.Text = " Menggunakan PHP IntelliSense"
}
$pos = $after2.End

$r3 = Find-From "3. " $pos
$r3.Text = "3."
$pos = $r3.End

# --- Change 4: remove the empty paragraph right after the second table,
#     just before the second "Saya Sudah Belajar dan Mengerti dan Saya Bisa" heading.
$r = Find-From "Saya Sudah Belajar dan Mengerti" $pos
$emptyPara = $d.Range($r.Start - 1, $r.Start - 1)
$emptyPara.Expand(4) | Out-Null
$emptyPara.Delete()
$pos = $r.End

# --- Change 5: update item 5 of the GitHub "Saya Sudah Belajar" list.
$r = Find-From "5. Mengerti istilah istilah pada GitHub." $pos
$r.Text = "5. Mengubah isi file yang tersimpan pada GitHub."
$pos = $r.End

# --- Change 6: fill in the second "Saya Belum Mengerti" list (items 1-2) and move
#     the _GoBack bookmark to sit after item 2.
$r = Find-From "Saya Belum Mengerti" $pos
$pos = $r.End

$r1 = Find-From "1." $pos
$ip1 = $d.Range($r1.End, $r1.End)
$ip1.InsertAfter(" Menghapus file yang terlanjur di masukkan ke repository")
$pos = $ip1.End

$r2 = Find-From "2." $pos
$ip2 = $d.Range($r2.End, $r2.End)
$ip2.InsertAfter(" ")
$bmRng = $d.Range($ip2.End, $ip2.End)
$d.Bookmarks.Add("_GoBack", $bmRng)
$pos = $ip2.End

# --- Change 7: remove the trailing empty paragraph at the very end of the document.
$endRng = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endRng.Expand(4) | Out-Null
if ($endRng.Text -eq "") {
    $endRng.Delete()
}
